$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Add the new "Description" column (E) ---

# Header cell: same style as the other header cells (C1/D1)
$ws.Range("E1").Value = "Description"
$ws.Range("D1").Copy()
$ws.Range("E1").PasteSpecial(-4122)  # xlPasteFormats (keeps the value already set)

# Data cells - written in the same order the original author filled them in
# (rows 4-8 first, then rows 2-3), so the shared-string table comes out in
# the same order as the target workbook.
$ws.Range("E4").Value = "The user needs to be able to see the state of the GVL. This task is trivial as a few team members have worked with displaying data to a screen in the past."
$ws.Range("E5").Value = "Battery data needs to be logged to nonvolatile storage, so the AER team can review it after a race. This risk is mitigated because there are many tutorials showing how to log data to an SD card using SPI."
$ws.Range("E6").Value = "The user needs to be able to request the logging of battery info. This task is trivial because it is easily accomplished by hooking up a momentary switch a GPIO pin on the MCU to toggle an interupt to begin logging data."
$ws.Range("E7").Value = "The user needs to bale to easily recharge the battery pack of the GLV device. This risk is mitigated, as we plan to use a laptop charger to supply external power to the device, recharging the batteries in a reasonable amount of time."
$ws.Range("E8").Value = "The team needs to layout a PCB for our design. No one on the team has any experience with this. The team has met with Mark Bruno in an effort to mitigate this risk."
$ws.Range("E2").Value = "The BMS IC is the heart of the project. This component is resposible for all aspects of charging, discharging, and battery protection features. If the team is not able to get the BMS IC working the entire project is dead. This risk has been mitigated by reviewing the datasheet for the IC, as common application circuits are given including supporting component values. "
$ws.Range("E3").Value = "The conversion from battery voltage to separate voltage rails of 24V and 12V to supply power to the low voltage electronics on the vehicle is the whole point of the project. If the team cannot make two voltage raills that can supply the current demanded that project is in trouble! This risk is mitigated as the team plans on using switching regulators to accomplish this task in an effiecent and eloquent manner."

# Widen the new column to fit its header/content, matching the other bestFit columns
$ws.Columns.Item(5).ColumnWidth = 14

# --- Update the view/selection state ---
$ws.Range("K10").Select() | Out-Null
